$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting the existing row 34 (and below) down to row 35.
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with the new data record.
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = Get-Date -Year 2022 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Range("D34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = 100112010
$ws.Range("G34").Value = "Achicoria"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 160
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = 9500
$ws.Range("N34").Value = "$/caja 18 unidades"
$ws.Range("O34").Value = "Región Metropolitana"
$ws.Range("P34").Value = 528
$ws.Range("Q34").Value = 18
$ws.Range("R34").Value = "Hortaliza"
